$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($r = 2; $r -le 279; $r++) {
    $ws.Cells.Item($r, 3).Value = 46075
}
